$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- sheet2 (kategorikas_listjumlahkas): "jenis" enum column becomes "jenis_id" FK column ---
$ws2.Range("C1").Value = "jenis_id"
$ws2.Range("C2").Value = -1
$ws2.Range("C3").Value = -2
$ws2.Columns("C").ColumnWidth = 7.17

# --- sheet3 (transaksikas): "jenis" enum column becomes "jenis_id" FK column ---
$ws3.Range("E1").Value = "jenis_id"
$ws3.Range("E2").Value = -1
$ws3.Range("E3").Value = -2
$ws3.Columns("E").ColumnWidth = 7.17

# --- add the new "jenistransaksikas" entity sheet after transaksikas ---
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws3)
$ws4.Name = "jenistransaksikas"

# Header row - reuse the bold/gray header style (s=3) from the kategorikas sheet
$ws1.Range("A1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)
$ws4.Range("A1").Value = "id"

$ws1.Range("B1").Copy()
$ws4.Range("B1").PasteSpecial(-4122)
$ws4.Range("B1").Value = "createdBy"

$ws1.Range("C1").Copy()
$ws4.Range("C1").PasteSpecial(-4122)
$ws4.Range("C1").Value = "createdDate"

$ws1.Range("D1").Copy()
$ws4.Range("D1").PasteSpecial(-4122)
$ws4.Range("D1").Value = "deleted"

$ws1.Range("F1").Copy()
$ws4.Range("E1").PasteSpecial(-4122)
$ws4.Range("E1").Value = "modifiedBy"

$ws1.Range("G1").Copy()
$ws4.Range("F1").PasteSpecial(-4122)
$ws4.Range("F1").Value = "modifiedDate"

$ws1.Range("H1").Copy()
$ws4.Range("G1").PasteSpecial(-4122)
$ws4.Range("G1").Value = "nama"

# Data row 2
$ws4.Range("A2").Value = -1
$ws4.Range("B2").Value = "admin"
$ws1.Range("C3").Copy()
$ws4.Range("C2").PasteSpecial(-4122)
$ws4.Range("C2").Value = 41971
$ws4.Range("G2").Value = "Dalam Kota"

# Data row 3
$ws4.Range("A3").Value = -2
$ws4.Range("B3").Value = "admin"
$ws1.Range("C3").Copy()
$ws4.Range("C3").PasteSpecial(-4122)
$ws4.Range("C3").Value = 41971
$ws4.Range("G3").Value = "Luar Kota"

# Column widths for the new sheet (closest achievable to the authored widths)
$ws4.Columns("A").ColumnWidth = 1.83
$ws4.Columns("B").ColumnWidth = 9.33
$ws4.Columns("C").ColumnWidth = 11.17
$ws4.Columns("D").ColumnWidth = 7
$ws4.Columns("E").ColumnWidth = 10.5
$ws4.Columns("F").ColumnWidth = 12.5
$ws4.Columns("G").ColumnWidth = 9.83

# --- restore / update cursor + active-tab state on each sheet ---
$ws1.Select()
$ws1.Range("I1").Select()

$ws2.Select()
$ws2.Range("C3").Select()

$ws3.Select()
$ws3.Range("E2").Select()
$excel.ActiveWindow.ScrollColumn = 3

$ws4.Select()
$ws4.Range("C2").Select()

Write-Host "done"
